$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(86, "CDMO_PRD_AHM", "Synology", "NAS", "2025-12-25 23:00:50", "2025-12-25 23:08:32"),
    @(86, "supervisor", "Fortinet", "FortiSIEM", "2025-12-25 23:00:50", "2025-12-25 23:00:50"),
    @(86, "MEH_PRD_NAS", "Synology", "NAS", "2025-12-25 23:00:52", "2025-12-25 23:08:29"),
    @(86, "LAMFILE_NAS", "Synology", "NAS", "2025-12-25 23:00:53", "2025-12-25 23:08:26"),
    @(86, "ARCHIVAL_AHM_NAS", "Synology", "NAS", "2025-12-25 23:00:53", "2025-12-25 23:08:29"),
    @(86, "MEHLAMFILE_NAS", "Synology", "NAS", "2025-12-25 23:00:53", "2025-12-25 23:08:32"),
    @(86, "Lambda-Changodar-60F", "Fortinet", "Fortigate", "2025-12-25 23:00:56", "2025-12-25 23:00:56"),
    @(86, "Fortigate_Firewall_Ahm_OT", "Fortinet", "Fortigate", "2025-12-25 23:00:57", "2025-12-25 23:08:29"),
    @(86, "Lambda_MEH_FW-100F", "Fortinet", "Fortigate", "2025-12-25 23:00:59", "2025-12-25 23:08:29"),
    @(86, "AHM_PRD_NAS", "Synology", "NAS", "2025-12-25 23:05:18", "2025-12-25 23:08:26"),
    @(86, "HOST-3.123.68.65", "Fortinet", "FortiWeb Cloud", "2025-12-25 23:05:22", "2025-12-25 23:08:23"),
    @(86, "LTRADC.ahm.lambdacro.com", "Microsoft", "Windows", "2025-12-25 23:07:47", "2025-12-25 23:07:47"),
    @(86, "LTRPDC.ahm.lambdacro.com", "Microsoft", "Windows", "2025-12-25 23:07:52", "2025-12-25 23:07:52"),
    @(86, "Lambda.coll", "ManageEngine", "ADAuditPlus", "2025-12-25 23:08:06", "2025-12-25 23:08:06"),
    @(86, "PV_RPT_DB.ahm.lambdacro.com", "Microsoft", "Windows Server 2022", "2025-12-25 23:08:26", "2025-12-25 23:08:26"),
    @(86, "HOST-192.168.22.98", "Fortinet", "FortiSIEM", "2025-12-25 23:08:26", "2025-12-25 23:08:26"),
    @(86, "BIZNET7.ahm.lambdacro.com", "Microsoft", "Windows Server 2022", "2025-12-25 23:08:29", "2025-12-25 23:08:29"),
    @(86, "AHM_VCenter", "VMware", "vCenter", "2025-12-25 23:08:42", "2025-12-25 23:08:42"),
    @(86, "Lambda_Primary_Firewall-200F", "Fortinet", "Fortigate", "2025-12-25 23:08:53", "2025-12-25 23:08:53"),
    @(86, "HOST-192.168.22.98,Lambda.coll", "Fortinet", "FortiSIEM", "2025-12-25 23:09:46", "2025-12-25 23:09:46")
)

$rowIndex = 2
foreach ($row in $data) {
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $row[$col - 1]
    }
    $rowIndex++
}
